$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.284.53"
$ws.Range("E2").Value = "  -2.71%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.215.87"
$ws.Range("E3").Value = "  -2.47%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "107.90"
$ws.Range("E5").Value = "  -11.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "296.12"
$ws.Range("E6").Value = "  +11.51%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -3.21%  "
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.594"
$ws.Range("E9").Value = "  -4.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "43.43"
$ws.Range("E10").Value = "  -9.21%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0909"
$ws.Range("E11").Value = "  -3.72%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.49"
$ws.Range("E12").Value = "  +0.61%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.77"
$ws.Range("E13").Value = "  -5.80%  "
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.962"
$ws.Range("E15").Value = "  +5.49%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.96"
$ws.Range("E16").Value = "  -3.01%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.544.85"
$ws.Range("E17").Value = "  -2.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.225.23"
$ws.Range("E18").Value = "  -1.78%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "42.235.93"
$ws.Range("E19").Value = "  -3.00%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.37"
$ws.Range("E20").Value = "  +6.95%  "
$ws.Range("E21").Value = "  -4.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.40"
$ws.Range("E22").Value = "  +0.26%  "
$ws.Range("E23").Value = "  +20.70%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.30"
$ws.Range("E24").Value = "  -3.85%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "228.00"
$ws.Range("E25").Value = "  -3.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.98"
$ws.Range("E26").Value = "  -5.77%  "
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "11.58"
$ws.Range("E28").Value = "  -3.06%  "
$ws.Range("E29").Value = "  -1.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.38"
$ws.Range("E30").Value = "  -8.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.19"
$ws.Range("E31").Value = "  -5.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "173.46"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.91"
$ws.Range("E33").Value = "  -3.41%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0889"
$ws.Range("E34").Value = "  -2.89%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.04"
$ws.Range("E35").Value = "  +9.63%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.52"
$ws.Range("E36").Value = "  -3.72%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.32"
$ws.Range("E37").Value = "  +1.03%  "
$ws.Range("E38").Value = "  -3.67%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0365"
$ws.Range("E39").Value = "  -3.05%  "
$ws.Range("E40").Value = "  -4.14%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.43"
$ws.Range("E41").Value = "  -4.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "71.75"
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.230"
$ws.Range("E43").Value = "  -3.04%  "
$ws.Range("E44").Value = "  +0.11%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.55"
$ws.Range("E45").Value = "  -9.57%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.30"
$ws.Range("E46").Value = "  -5.49%  "
$ws.Range("E47").Value = "  -5.81%  "
$ws.Range("E48").Value = "  +3.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.36"
$ws.Range("E49").Value = "  +1.90%  "
$ws.Range("B50").Value = "FraxShare"
$ws.Range("C50").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.39"
$ws.Range("E50").Value = "  -1.39%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.62"
$ws.Range("E51").Value = "  +3.75%  "
